$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "2022-Q1" sheet between "2021-Q4" and "总计" ---
$template = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Add($null, $template)
$ws.Name = "2022-Q1"

# Reuse the header-row / index-column formatting from the "2021-Q4" template
$template.Range("B1:H1").Copy($ws.Range("B1:H1"))
$template.Range("A2").Copy($ws.Range("A2:A49"))

# --- 2. Header text for the new sheet ---
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# --- 3. Fund-level holdings data (rows 2-49) ---
# row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "161810"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "银华内需精选混合(LOF)"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "94.71"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "8.85"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2.2647"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 1
# row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "009394"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "银华同力精选混合"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "20.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "94.68"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "8.27"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "1.6565"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = 2
# row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "240022"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "华宝资源优选混合A"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "25.75"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "85.32"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "2.74"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "0.7056"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = 10
# row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "162207"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "泰达宏利效率优选混合(LOF)"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "6.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "69.26"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "6.32"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "0.3925"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 2
# row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "001170"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "泰达宏利复兴伟业灵活配置混合"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "92.00"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "8.23"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "0.2527"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = 2
# row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "005273"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "华商可转债债券A"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "10.06"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "39.60"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "2.51"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "0.2525"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = 4
# row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "011068"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "华宝资源优选混合C"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "9.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "85.32"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "2.74"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "0.2502"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = 10
# row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "012138"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "景顺长城安益回报一年持有期混合型证券投资基金A"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "32.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "25.23"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "0.75"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "0.2473"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = 8
# row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "510410"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "博时上证自然资源ETF"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "4.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "98.61"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "4.44"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "0.2073"
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value = 5
# row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "159930"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "汇添富中证能源ETF"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "2.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "99.19"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "9.00"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "0.1908"
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").Value = 4
# row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "470007"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "汇添富上证综合指数"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "94.43"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "2.33"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0.1754"
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value = 5
# row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "519767"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "交银施罗德科技创新灵活配置混合"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "90.38"
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "3.96"
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "0.1655"
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").Value = 5
# row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "510210"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "富国上证综指ETF"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "99.46"
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "2.44"
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0.1654"
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").Value = 4
# row 15
$ws.Range("A15").Value = 13
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "217012"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "招商行业领先混合A"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "80.77"
$ws.Range("E15").Style = "Normal"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "4.91"
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0.1404"
$ws.Range("G15").Style = "Normal"
$ws.Range("H15").Value = 8
# row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "960019"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "招商行业领先混合H"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "80.77"
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "4.91"
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "0.1404"
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").Value = 8
# row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "005284"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "华商可转债债券C"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "39.60"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "2.51"
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "0.0954"
$ws.Range("G17").Style = "Normal"
$ws.Range("H17").Value = 4
# row 18
$ws.Range("A18").Value = 16
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "161217"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "国投瑞银中证上游资源产业指数(LOF)"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "94.10"
$ws.Range("E18").Style = "Normal"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "3.00"
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "0.0930"
$ws.Range("G18").Style = "Normal"
$ws.Range("H18").Value = 9
# row 19
$ws.Range("A19").Value = 17
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "009141"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "泰达宏利价值长青混合A"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "90.52"
$ws.Range("E19").Style = "Normal"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "4.08"
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "0.0910"
$ws.Range("G19").Style = "Normal"
$ws.Range("H19").Value = 9
# row 20
$ws.Range("A20").Value = 18
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "011018"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "景顺长城安泽回报一年持有期混合A"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "34.82"
$ws.Range("E20").Style = "Normal"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "1.00"
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "0.0781"
$ws.Range("G20").Style = "Normal"
$ws.Range("H20").Value = 5
# row 21
$ws.Range("A21").Value = 19
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "005317"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "万家瑞舜灵活配置混合A"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "22.90"
$ws.Range("E21").Style = "Normal"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "0.75"
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "0.0723"
$ws.Range("G21").Style = "Normal"
$ws.Range("H21").Value = 6
# row 22
$ws.Range("A22").Value = 20
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "540002"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "汇丰晋信龙腾混合"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "66.68"
$ws.Range("E22").Style = "Normal"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "2.00"
$ws.Range("F22").Style = "Normal"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0.0658"
$ws.Range("G22").Style = "Normal"
$ws.Range("H22").Value = 10
# row 23
$ws.Range("A23").Value = 21
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "690008"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "民生加银中证内地资源主题指数"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "94.56"
$ws.Range("E23").Style = "Normal"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "3.12"
$ws.Range("F23").Style = "Normal"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "0.0583"
$ws.Range("G23").Style = "Normal"
$ws.Range("H23").Value = 9
# row 24
$ws.Range("A24").Value = 22
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "000368"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "汇添富沪深300安中指数"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "94.41"
$ws.Range("E24").Style = "Normal"
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value = "2.49"
$ws.Range("F24").Style = "Normal"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "0.0560"
$ws.Range("G24").Style = "Normal"
$ws.Range("H24").Value = 7
# row 25
$ws.Range("A25").Value = 23
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "001678"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "英大国企改革主题股票"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "93.10"
$ws.Range("E25").Style = "Normal"
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = "7.56"
$ws.Range("F25").Style = "Normal"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "0.0552"
$ws.Range("G25").Style = "Normal"
$ws.Range("H25").Value = 7
# row 26
$ws.Range("A26").Value = 24
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "004731"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "万家瑞尧灵活配置混合A"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "23.06"
$ws.Range("E26").Style = "Normal"
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "0.75"
$ws.Range("F26").Style = "Normal"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "0.0541"
$ws.Range("G26").Style = "Normal"
$ws.Range("H26").Value = 5
# row 27
$ws.Range("A27").Value = 25
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "001635"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "万家瑞益灵活配置混合A"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "23.91"
$ws.Range("E27").Style = "Normal"
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0.74"
$ws.Range("F27").Style = "Normal"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "0.0519"
$ws.Range("G27").Style = "Normal"
$ws.Range("H27").Value = 6
# row 28
$ws.Range("A28").Value = 26
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "510170"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "国联安上证大宗商品股票ETF"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "98.14"
$ws.Range("E28").Style = "Normal"
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "2.24"
$ws.Range("F28").Style = "Normal"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "0.0497"
$ws.Range("G28").Style = "Normal"
$ws.Range("H28").Value = 8
# row 29
$ws.Range("A29").Value = 27
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "001488"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "万家瑞丰灵活配置混合A"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "20.31"
$ws.Range("E29").Style = "Normal"
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0.64"
$ws.Range("F29").Style = "Normal"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "0.0381"
$ws.Range("G29").Style = "Normal"
$ws.Range("H29").Value = 5
# row 30
$ws.Range("A30").Value = 28
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "159945"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "广发中证全指能源ETF"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "98.90"
$ws.Range("E30").Style = "Normal"
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "7.05"
$ws.Range("F30").Style = "Normal"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "0.0374"
$ws.Range("G30").Style = "Normal"
$ws.Range("H30").Value = 4
# row 31
$ws.Range("A31").Value = 29
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "001636"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "万家瑞益灵活配置混合C"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "23.91"
$ws.Range("E31").Style = "Normal"
$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = "0.74"
$ws.Range("F31").Style = "Normal"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0.0332"
$ws.Range("G31").Style = "Normal"
$ws.Range("H31").Value = 6
# row 32
$ws.Range("A32").Value = 30
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "011534"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "万家民瑞祥明6个月持有期混合型证券投资基金A"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "23.07"
$ws.Range("E32").Style = "Normal"
$ws.Range("F32").NumberFormat = "@"
$ws.Range("F32").Value = "0.72"
$ws.Range("F32").Style = "Normal"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "0.0306"
$ws.Range("G32").Style = "Normal"
$ws.Range("H32").Value = 9
# row 33
$ws.Range("A33").Value = 31
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "011097"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "达诚宜创精选混合A"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "89.58"
$ws.Range("E33").Style = "Normal"
$ws.Range("F33").NumberFormat = "@"
$ws.Range("F33").Value = "2.06"
$ws.Range("F33").Style = "Normal"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "0.0220"
$ws.Range("G33").Style = "Normal"
$ws.Range("H33").Value = 7
# row 34
$ws.Range("A34").Value = 32
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "516570"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "易方达中证石化产业交易型开放式指数证券投资基金"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "96.03"
$ws.Range("E34").Style = "Normal"
$ws.Range("F34").NumberFormat = "@"
$ws.Range("F34").Value = "5.55"
$ws.Range("F34").Style = "Normal"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "0.0200"
$ws.Range("G34").Style = "Normal"
$ws.Range("H34").Value = 5
# row 35
$ws.Range("A35").Value = 33
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "004335"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "华宝新飞跃灵活配置混合"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "44.86"
$ws.Range("E35").Style = "Normal"
$ws.Range("F35").NumberFormat = "@"
$ws.Range("F35").Value = "0.65"
$ws.Range("F35").Style = "Normal"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "0.0179"
$ws.Range("G35").Style = "Normal"
$ws.Range("H35").Value = 10
# row 36
$ws.Range("A36").Value = 34
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "004732"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "万家瑞尧灵活配置混合C"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "23.06"
$ws.Range("E36").Style = "Normal"
$ws.Range("F36").NumberFormat = "@"
$ws.Range("F36").Value = "0.75"
$ws.Range("F36").Style = "Normal"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "0.0176"
$ws.Range("G36").Style = "Normal"
$ws.Range("H36").Value = 5
# row 37
$ws.Range("A37").Value = 35
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "005328"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "前海开源价值策略股票"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "92.34"
$ws.Range("E37").Style = "Normal"
$ws.Range("F37").NumberFormat = "@"
$ws.Range("F37").Value = "4.00"
$ws.Range("F37").Style = "Normal"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "0.0148"
$ws.Range("G37").Style = "Normal"
$ws.Range("H37").Value = 6
# row 38
$ws.Range("A38").Value = 36
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "001489"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "万家瑞丰灵活配置混合C"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "20.31"
$ws.Range("E38").Style = "Normal"
$ws.Range("F38").NumberFormat = "@"
$ws.Range("F38").Value = "0.64"
$ws.Range("F38").Style = "Normal"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "0.0148"
$ws.Range("G38").Style = "Normal"
$ws.Range("H38").Value = 5
# row 39
$ws.Range("A39").Value = 37
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "001789"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "国泰量化收益灵活配置混合"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "91.40"
$ws.Range("E39").Style = "Normal"
$ws.Range("F39").NumberFormat = "@"
$ws.Range("F39").Value = "2.43"
$ws.Range("F39").Style = "Normal"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "0.0139"
$ws.Range("G39").Style = "Normal"
$ws.Range("H39").Value = 9
# row 40
$ws.Range("A40").Value = 38
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "161816"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "银华中证等权重90指数（LOF）"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "92.74"
$ws.Range("E40").Style = "Normal"
$ws.Range("F40").NumberFormat = "@"
$ws.Range("F40").Value = "1.34"
$ws.Range("F40").Style = "Normal"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "0.0118"
$ws.Range("G40").Style = "Normal"
$ws.Range("H40").Value = 5
# row 41
$ws.Range("A41").Value = 39
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "012139"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "景顺长城安益回报一年持有期混合型证券投资基金C"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "25.23"
$ws.Range("E41").Style = "Normal"
$ws.Range("F41").NumberFormat = "@"
$ws.Range("F41").Value = "0.75"
$ws.Range("F41").Style = "Normal"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "0.0103"
$ws.Range("G41").Style = "Normal"
$ws.Range("H41").Value = 8
# row 42
$ws.Range("A42").Value = 40
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "510190"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "华安上证龙头ETF"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "97.53"
$ws.Range("E42").Style = "Normal"
$ws.Range("F42").NumberFormat = "@"
$ws.Range("F42").Value = "1.15"
$ws.Range("F42").Style = "Normal"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "0.0061"
$ws.Range("G42").Style = "Normal"
$ws.Range("H42").Value = 9
# row 43
$ws.Range("A43").Value = 41
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "009142"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "泰达宏利价值长青混合C"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "90.52"
$ws.Range("E43").Style = "Normal"
$ws.Range("F43").NumberFormat = "@"
$ws.Range("F43").Value = "4.08"
$ws.Range("F43").Style = "Normal"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "0.0045"
$ws.Range("G43").Style = "Normal"
$ws.Range("H43").Value = 9
# row 44
$ws.Range("A44").Value = 42
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "011098"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "达诚宜创精选混合C"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "89.58"
$ws.Range("E44").Style = "Normal"
$ws.Range("F44").NumberFormat = "@"
$ws.Range("F44").Value = "2.06"
$ws.Range("F44").Style = "Normal"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "0.0041"
$ws.Range("G44").Style = "Normal"
$ws.Range("H44").Value = 7
# row 45
$ws.Range("A45").Value = 43
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "011535"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "万家民瑞祥明6个月持有期混合型证券投资基金C"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "23.07"
$ws.Range("E45").Style = "Normal"
$ws.Range("F45").NumberFormat = "@"
$ws.Range("F45").Value = "0.72"
$ws.Range("F45").Style = "Normal"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "0.0033"
$ws.Range("G45").Style = "Normal"
$ws.Range("H45").Value = 9
# row 46
$ws.Range("A46").Value = 44
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "005960"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "博时量化价值股票A"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "90.69"
$ws.Range("E46").Style = "Normal"
$ws.Range("F46").NumberFormat = "@"
$ws.Range("F46").Value = "0.96"
$ws.Range("F46").Style = "Normal"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "0.0029"
$ws.Range("G46").Style = "Normal"
$ws.Range("H46").Value = 9
# row 47
$ws.Range("A47").Value = 45
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "011019"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "景顺长城安泽回报一年持有期混合C"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "34.82"
$ws.Range("E47").Style = "Normal"
$ws.Range("F47").NumberFormat = "@"
$ws.Range("F47").Value = "1.00"
$ws.Range("F47").Style = "Normal"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "0.0024"
$ws.Range("G47").Style = "Normal"
$ws.Range("H47").Value = 5
# row 48
$ws.Range("A48").Value = 46
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "005961"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "博时量化价值股票C"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "90.69"
$ws.Range("E48").Style = "Normal"
$ws.Range("F48").NumberFormat = "@"
$ws.Range("F48").Value = "0.96"
$ws.Range("F48").Style = "Normal"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "0.0023"
$ws.Range("G48").Style = "Normal"
$ws.Range("H48").Value = 9
# row 49
$ws.Range("A49").Value = 47
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "005318"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "万家瑞舜灵活配置混合C"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "22.90"
$ws.Range("E49").Style = "Normal"
$ws.Range("F49").NumberFormat = "@"
$ws.Range("F49").Value = "0.75"
$ws.Range("F49").Style = "Normal"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "0.0020"
$ws.Range("G49").Style = "Normal"
$ws.Range("H49").Value = 6

# --- 4. Update the "总计" (totals) summary sheet: insert a 2022-Q1 row ---
$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert()
$tot.Range("B2:D2").ClearFormats()
$tot.Range("A3").Copy($tot.Range("A2"))
$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 48
$tot.Range("D2").Value = 8.34

# Renumber the index column for the rows that shifted down
$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4
$tot.Range("A7").Value = 5

